$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7219566702842712
$ws.Range("B1").Value = 1.05582857131958
$ws.Range("C1").Value = 3.401817560195923
$ws.Range("D1").Value = 3.745327711105347
$ws.Range("E1").Value = 2.054075479507446
